# Add a "batch" / survival summary table to the right of the existing data
# (columns K:O, rows 7-18) on Sheet1, per commit "added batch to survival script".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (bold) ---------------------------------------------------
$ws.Range("K7").Value = "Treatment"
$ws.Range("L7").Value = "Lifestage"
$ws.Range("M7").Value = "Age"
$ws.Range("N7").Value = "Replicate-Tanks-per-Treatment"
$ws.Range("O7").Value = "Total-Samples"
$ws.Range("K7:O7").Font.Bold = $true

# --- Data rows -------------------------------------------------------------
# Row 8: Control / Larvae / 3 days / n=1 / 15
$ws.Range("K8").Value = "Control"
$ws.Range("L8").Value = "Larvae"
$ws.Range("M8").Value = "3 days"
$ws.Range("N8").Value = "n=1"
$ws.Range("O8").Value = 15

# Row 9: Treated / Larvae / 3 days / n=1 / 15
$ws.Range("K9").Value = "Treated"
$ws.Range("L9").Value = "Larvae"
$ws.Range("M9").Value = "3 days"
$ws.Range("N9").Value = "n=1"
$ws.Range("O9").Value = 15

# Row 10: Control / Larvae / 8 days / n=2-3 / 15
$ws.Range("K10").Value = "Control"
$ws.Range("L10").Value = "Larvae"
$ws.Range("M10").Value = "8 days"
$ws.Range("N10").Value = "n=2-3"
$ws.Range("O10").Value = 15

# Row 11: Treated / Larvae / 8 days / n=2-3 / 15
$ws.Range("K11").Value = "Treated"
$ws.Range("L11").Value = "Larvae"
$ws.Range("M11").Value = "8 days"
$ws.Range("N11").Value = "n=2-3"
$ws.Range("O11").Value = 15

# Row 12: Control / Larvae / 14 days / n=2-3 / 15
$ws.Range("K12").Value = "Control"
$ws.Range("L12").Value = "Larvae"
$ws.Range("M12").Value = "14 days"
$ws.Range("N12").Value = "n=2-3"
$ws.Range("O12").Value = 15

# Row 13: Treated / Larvae / 14 days / n=2-3 / 15
$ws.Range("K13").Value = "Treated"
$ws.Range("L13").Value = "Larvae"
$ws.Range("M13").Value = "14 days"
$ws.Range("N13").Value = "n=2-3"
$ws.Range("O13").Value = 15

# Row 14: Control / Spat / 60 days / n=3 / 15
$ws.Range("K14").Value = "Control"
$ws.Range("L14").Value = "Spat"
$ws.Range("M14").Value = "60 days"
$ws.Range("N14").Value = "n=3"
$ws.Range("O14").Value = 15

# Row 15: Treated / Spat / 60 days / n=3 / 15
$ws.Range("K15").Value = "Treated"
$ws.Range("L15").Value = "Spat"
$ws.Range("M15").Value = "60 days"
$ws.Range("N15").Value = "n=3"
$ws.Range("O15").Value = 15

# Row 16: Control / Seed / 135 days / n=2 / 20
$ws.Range("K16").Value = "Control"
$ws.Range("L16").Value = "Seed"
$ws.Range("M16").Value = "135 days"
$ws.Range("N16").Value = "n=2"
$ws.Range("O16").Value = 20

# Row 17: Treated / Seed / 135 days / n=2 / 20
$ws.Range("K17").Value = "Treated"
$ws.Range("L17").Value = "Seed"
$ws.Range("M17").Value = "135 days"
$ws.Range("N17").Value = "n=2"
$ws.Range("O17").Value = 20

# Row 18: TOTAL label + sum formula
$ws.Range("N18").Value = "TOTAL"
$ws.Range("O18").Formula = "=SUM(O8:O17)"

# --- Column widths for the new columns --------------------------------
# (ColumnWidth values chosen so the engine's internal pixel snapping lands
# on the stored widths 16.33203125 / 18.1640625 from the target file)
$ws.Range("N1").ColumnWidth = 15.42
$ws.Range("O1").ColumnWidth = 17.25

# --- View state (matches updated pane / selection in the diff) ---------
[void]$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
[void]$ws.Range("L22").Select()
